# Auto-generated Excel COM-interop edit script
# Applies the cell-level numeric updates described by the commit diff
# (market price refresh / profit recompute across the Leve profit tables).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 902.7857
$ws.Range("J33").Value = 2239.4
$ws.Range("L33").Value = 2239.4
$ws.Range("N33").Value = -2697.4
$ws.Range("H80").Value = 6450
$ws.Range("J80").Value = 6450
$ws.Range("L80").Value = 19350
$ws.Range("N80").Value = -21346
$ws.Range("H83").Value = 6450
$ws.Range("J83").Value = 6450
$ws.Range("L83").Value = 58050
$ws.Range("N83").Value = -68034
$ws.Range("H97").Value = 2285.6667
$ws.Range("J97").Value = 2285.6667
$ws.Range("L97").Value = 6857.000100000001
$ws.Range("N97").Value = -7849.000100000001
$ws.Range("H138").Value = 10512.467
$ws.Range("I138").Value = 11447.75
$ws.Range("J138").Value = 10172.363
$ws.Range("K138").Value = 34343.25
$ws.Range("L138").Value = 30517.089
$ws.Range("M138").Value = -29203.25
$ws.Range("N138").Value = -40797.089
$ws.Range("H141").Value = 2500.9
$ws.Range("I141").Value = 2378.889
$ws.Range("K141").Value = 7136.667
$ws.Range("M141").Value = -1956.667

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4000
$ws.Range("J74").Value = 4000
$ws.Range("L74").Value = 4000
$ws.Range("N74").Value = -5748
$ws.Range("H77").Value = 4000
$ws.Range("J77").Value = 4000
$ws.Range("L77").Value = 20000
$ws.Range("N77").Value = -28736
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("N86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("N89").Value = 0
$ws.Range("H102").Value = 2000
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value = 3351
$ws.Range("I110").Value = 3351
$ws.Range("K110").Value = 3351
$ws.Range("M110").Value = -1306
$ws.Range("H119").Value = 33999.5
$ws.Range("J119").Value = 33999.5
$ws.Range("L119").Value = 33999.5
$ws.Range("N119").Value = -43675.5
$ws.Range("H122").Value = 1208.6666
$ws.Range("I122").Value = 1273.091
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 3819.273
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = -1369.273
$ws.Range("N122").Value = -6400

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 8185.5
$ws.Range("I26").Value = 8185.5
$ws.Range("K26").Value = 8185.5
$ws.Range("M26").Value = -7893.5
$ws.Range("H86").Value = 4104.647
$ws.Range("I86").Value = 4270.5625
$ws.Range("J86").Value = 1450
$ws.Range("K86").Value = 4270.5625
$ws.Range("L86").Value = 1450
$ws.Range("M86").Value = -3147.5625
$ws.Range("N86").Value = -3696
$ws.Range("H89").Value = 4104.647
$ws.Range("I89").Value = 4270.5625
$ws.Range("J89").Value = 1450
$ws.Range("K89").Value = 21352.8125
$ws.Range("L89").Value = 7250
$ws.Range("M89").Value = -15736.8125
$ws.Range("N89").Value = -18482
$ws.Range("H99").Value = 2437.5
$ws.Range("I99").Value = 2437.5
$ws.Range("K99").Value = 2437.5
$ws.Range("M99").Value = -939.5
$ws.Range("H134").Value = 1221.25
$ws.Range("I134").Value = 1221.25
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3663.75
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -1128.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = 0
$ws.Range("H105").Value = 10602.4
$ws.Range("I105").Value = 9006
$ws.Range("K105").Value = 9006
$ws.Range("M105").Value = -7259
$ws.Range("H132").Value = 3706.75
$ws.Range("I132").Value = 3252.6924
$ws.Range("K132").Value = 9758.0772
$ws.Range("M132").Value = -7228.0772
$ws.Range("H134").Value = 1347.9231
$ws.Range("I134").Value = 1356.0834
$ws.Range("K134").Value = 4068.2502
$ws.Range("M134").Value = -1533.2502

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 55.875
$ws.Range("I2").Value = 20
$ws.Range("J2").Value = 115.666664
$ws.Range("K2").Value = 120
$ws.Range("L2").Value = 693.999984
$ws.Range("M2").Value = -7
$ws.Range("N2").Value = -919.999984
$ws.Range("H3").Value = 5000
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 15000
$ws.Range("N3").Value = -15224
$ws.Range("H15").Value = 1916.6666
$ws.Range("I15").Value = 250
$ws.Range("J15").Value = 2750
$ws.Range("K15").Value = 750
$ws.Range("L15").Value = 8250
$ws.Range("M15").Value = -610
$ws.Range("N15").Value = -8530
$ws.Range("H17").Value = 7687.5
$ws.Range("I17").Value = 7687.5
$ws.Range("K17").Value = 23062.5
$ws.Range("M17").Value = -22893.5
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = 0

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").ClearContents()
$ws.Range("N38").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("N43").Value = 0
$ws.Range("H104").Value = 28000
$ws.Range("J104").Value = 28000
$ws.Range("L104").Value = 28000
$ws.Range("N104").Value = -34988

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4166.6665
$ws.Range("I81").Value = 4166.6665
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 8333.333000000001
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -7272.333000000001
$ws.Range("H84").Value = 4166.6665
$ws.Range("I84").Value = 4166.6665
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 41666.665
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -36362.665
$ws.Range("H100").Value = 3873250.8
$ws.Range("I100").Value = 4979092
$ws.Range("J100").Value = 2806
$ws.Range("K100").Value = 9958184
$ws.Range("L100").Value = 5612
$ws.Range("M100").Value = -9957643
$ws.Range("N100").Value = -6694
$ws.Range("H119").Value = 38194
$ws.Range("J119").Value = 38194
$ws.Range("L119").Value = 38194
$ws.Range("N119").Value = -47870
$ws.Range("H122").Value = 2780.75
$ws.Range("I122").Value = 2780.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8342.25
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -5892.25

Write-Output "Applied Sheets profit/price updates across ALC, ARM, BSM, CRP, CUL, LTW, WVR"